$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Row 32 - add a SATURDAY VL credit line (value 1.25 in EARNED)
# ---------------------------------------------------------------------------
$ws.Range("C32").Value = 1.25

# ---------------------------------------------------------------------------
# Row 33 - new SL(1-0-0) entry, 1.25 earned, 1 day absence, paid through 10/25/2023
# ---------------------------------------------------------------------------
$ws.Range("B33").Value = "SL(1-0-0)"
$ws.Range("C33").Value = 1.25
$ws.Range("H33").Value = 1
# give K33 the same date format (style) already used by the other "paid-through" cells
$ws.Range("K32").Copy()
$ws.Range("K33").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K33").Value = 45224

# ---------------------------------------------------------------------------
# Row 34 - November period starts, VL(1-0-0) entry, paid through 11/22/2023
# ---------------------------------------------------------------------------
$ws.Range("A34").Value = 45231
$ws.Range("B34").Value = "VL(1-0-0)"
$ws.Range("D34").Value = 1
$ws.Range("K32").Copy()
$ws.Range("K34").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K34").Value = 45252

# ---------------------------------------------------------------------------
# Rows 35-48 - fill in the remaining monthly PERIOD dates (Dec 2023 .. Dec 2024)
# Row 36 is the "2024" year banner (text, bold, quote-prefixed) like rows 10/23
# ---------------------------------------------------------------------------
$ws.Range("A35").Value = 45261

$ws.Range("A23").Copy()
$ws.Range("A36").PasteSpecial(-4122)   # xlPasteFormats (bold/quote-prefixed year style)
$ws.Range("A36").Value = "'2024"

$ws.Range("A37").Value = 45292
$ws.Range("A38").Value = 45323
$ws.Range("A39").Value = 45352
$ws.Range("A40").Value = 45383
$ws.Range("A41").Value = 45413
$ws.Range("A42").Value = 45444
$ws.Range("A43").Value = 45474
$ws.Range("A44").Value = 45505
$ws.Range("A45").Value = 45536
$ws.Range("A46").Value = 45566
$ws.Range("A47").Value = 45597
$ws.Range("A48").Value = 45627

# ---------------------------------------------------------------------------
# Extend the leave table by one row: old row 131 (the visually-different last
# row of the table) moves down to row 132, and the new row 131 takes on the
# regular interior-row formatting (copied from row 130).
# ---------------------------------------------------------------------------
$ws.Range("A131:K131").Copy()
$ws.Range("A132:K132").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A130:K130").Copy()
$ws.Range("A131:K131").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G132").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K132"))

# ---------------------------------------------------------------------------
# Update the saved view state: scroll the frozen-pane section down and move
# the selection in the bottom pane to K34.
# ---------------------------------------------------------------------------
$ws.Range("K34").Select()
